$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new, more recent weekly price record was added to the data set.
# It belongs right after the header block of existing "Granada" rows,
# at row 7 — so insert a fresh row there, which pushes the former
# rows 7-12 down to rows 8-13 (dimension grows from T12 to T13).
$ws.Rows(7).Insert()

# Populate the newly inserted row 7 with the new record's data.
$ws.Cells.Item(7, 1).Value  = 6
$ws.Cells.Item(7, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(7, 3).Value  = "Metropolitana"
$ws.Cells.Item(7, 4).Value  = 44721
$ws.Cells.Item(7, 5).Value  = 13
$ws.Cells.Item(7, 6).Value  = "Fruta"
$ws.Cells.Item(7, 7).Value  = 100104
$ws.Cells.Item(7, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(7, 9).Value  = 100104001
$ws.Cells.Item(7, 10).Value = "Granada"
$ws.Cells.Item(7, 11).Value = "Wonderfull"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 7
$ws.Cells.Item(7, 14).Value = 300000
$ws.Cells.Item(7, 15).Value = 300000
$ws.Cells.Item(7, 16).Value = 300000
$ws.Cells.Item(7, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(7, 18).Value = "Región Metropolitana"
$ws.Cells.Item(7, 19).Value = 750
$ws.Cells.Item(7, 20).Value = 400
